$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: paragraph 2 ("Täida alloleva tabeli ...") gets a new
# lead-in sentence prepended and is split into two runs:
#   "Tühjenda allolev tabel ning t" + "äida alloleva tabeli ..."
# -----------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2Start = $p2.Range.Start
# Replace the leading "T" with the new prefix text (ends in lowercase "t")
$firstCharRange = $d.Range($p2Start, $p2Start + 1)
$prefix1 = "Tühjenda allolev tabel ning t"
$firstCharRange.Text = $prefix1
# Force a run split exactly at the boundary between the new prefix and
# the remaining original text by toggling a character formatting
# property on just the inserted span and reverting it.
$splitRange1 = $d.Range($p2Start, $p2Start + $prefix1.Length)
$splitRange1.Font.Bold = 1
$splitRange1.Font.Bold = 0

# -----------------------------------------------------------------
# Change 2 (part 1 of 2): the stray "_GoBack" bookmark that currently
# sits at the end of the "...100%." paragraph will be re-created at
# its new location further down (see Change 4); adding a bookmark
# with the same name moves it, so nothing else is required here.
# -----------------------------------------------------------------

# -----------------------------------------------------------------
# Change 3: paragraph "Salvesta pilt samasse kausta, kus on laevade
# pommitamise kood." is split into two runs:
#   "Salvesta " + "pilt kausta 'kood' alamkausta 'pildid'."
# -----------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p5Start = $p5.Range.Start
$p5End = $p5.Range.End
$lead = "Salvesta "
$tailNew = [char]0x2019 + "kood" + [char]0x2019 + " alamkausta " + [char]0x2019 + "pildid" + [char]0x2019
$newTail = "pilt kausta " + $tailNew + "."
$tailRange = $d.Range($p5Start + $lead.Length, $p5End - 1)
$tailRange.Text = $newTail
$splitRange2 = $d.Range($p5Start, $p5Start + $lead.Length)
$splitRange2.Font.Bold = 1
$splitRange2.Font.Bold = 0

# -----------------------------------------------------------------
# Change 4: a brand-new list paragraph is appended after "Vali
# laevade pommitamist mängides ..." (currently paragraph 6), one
# list-level deeper, holding two runs of text and ending with the
# relocated "_GoBack" bookmark.
# -----------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$p6Range = $p6.Range
$newParaPos = $p6Range.End
$p6Range.Collapse(0)
$p6Range.InsertParagraphAfter()

$p7 = $d.Paragraphs(7)
$p7.Range.ListFormat.ListLevelNumber = 2

$run3a = "Allolev t" + [char]0xE4 + "idetud tabel on pildi kujul failina " + [char]0x2019 + "naidis_laud.png" + [char]0x2019
$run3b = " ning ka seda saab m" + [char]0xE4 + "ngimiseks kasutada."
$placeholder = [char]0x7E + [char]0x7E + [char]0x7E + "PLACEHOLDER" + [char]0x7E + [char]0x7E + [char]0x7E

$p7 = $d.Paragraphs(7)
$p7.Range.InsertBefore($run3a)

$afterRun3aPos = $newParaPos + $run3a.Length
$r = $d.Range($afterRun3aPos, $afterRun3aPos)
$r.InsertBefore($run3b + $placeholder)

# Split runs between run3a and run3b.
$boldRange3 = $d.Range($newParaPos, $newParaPos + $run3a.Length)
$boldRange3.Font.Bold = 1
$boldRange3.Font.Bold = 0

# Re-home the "_GoBack" bookmark at the true end of the new text
# (the placeholder keeps this position away from the paragraph's
# final character slot while the bookmark is created).
$bmPos = $newParaPos + $run3a.Length + $run3b.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Now remove the placeholder text that trailed the bookmark.
$delRange = $d.Range($bmPos, $bmPos + $placeholder.Length)
$delRange.Text = ""
